# Added currentLichessUser to Frontend
#
# Appends a new challenge entry as row 3 on the (only) worksheet,
# mirroring the fields already recorded for row 2: same challenger,
# rating, wager, "accepted?" state and accepter placeholder, but for a
# newly created lichess game/challenge (new gameID, link and escrowID).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 3
$newRecord = @{
    1 = "mna6nw74"                         # gameID
    2 = "trashboatsr"                      # challenger
    3 = 1890                               # rating
    4 = 20                                 # wager
    5 = "https://lichess.org/mna6nw74"     # link
    6 = 4041                               # escrowID
    7 = $false                             # accepted?
    8 = "blank"                            # accepter
}

foreach ($col in $newRecord.Keys) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = $newRecord[$col]
    # Keep the new row's formatting consistent with the existing data row
    # (plain/default cell style) instead of silently inheriting whatever
    # display style the column itself carries.
    $cell.Style = "Normal"
}
